# Automatic update of files.
# Rows 3, 4 and 5 of the "Artfynd" sheet get cyclically rotated:
#   old row 4 -> row 3
#   old row 5 -> row 4
#   old row 3 -> row 5
# (columns C, K, P, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY are identical
#  across these rows / unchanged by the diff, so only the "record" columns
#  below are touched)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

function Set-RowValues($wsArg, $row, $values) {
    $wsArg.Range("A$row").Value = $values.A
    $wsArg.Range("B$row").Value = $values.B
    $wsArg.Range("D$row").Value = $values.D
    $wsArg.Range("E$row").Value = $values.E
    $wsArg.Range("F$row").Value = $values.F
    $wsArg.Range("G$row").Value = $values.G
    $wsArg.Range("H$row").Value = $values.H
    # Only touch column I when it actually needs a non-empty value; writing an
    # empty string removes the cell instead of leaving an empty placeholder.
    # The leading apostrophe forces the numeric-looking "35" to be stored as
    # text, matching the source data (Antal is a text column here).
    if ($null -ne $values.I) {
        $wsArg.Range("I$row").Value = "'" + $values.I
    }
    $wsArg.Range("Q$row").Value = $values.Q
    $wsArg.Range("R$row").Value = $values.R
    $wsArg.Range("S$row").Value = $values.S
    $wsArg.Range("Z$row").Value = $values.Z
    $wsArg.Range("AB$row").Value = $values.AB
}

# Snapshot of the original (pre-edit) row 3, 4 and 5 "record" values.
$origRow3 = [PSCustomObject]@{
    A = 111645826; B = 94134; D = "NT"; E = 53
    F = "Vedtrappmossa"; G = "Crossocalyx hellerianus"; H = "(Nees ex Lindenb.) Meyl."
    I = $null
    Q = 369469.7018694163; R = 6635346.076433333; S = 10
    Z = "13:25"; AB = "13:25"
}

$origRow4 = [PSCustomObject]@{
    A = 111644287; B = 96348; D = "VU"; E = 220787
    F = "Knärot"; G = "Goodyera repens"; H = "(L.) R. Br."
    I = "35"
    Q = 369410.2171064656; R = 6635288.297872287; S = 10
    Z = "12:37"; AB = "12:37"
}

$origRow5 = [PSCustomObject]@{
    A = 111644923; B = 56543; D = "NT"; E = 103021
    F = "Talltita"; G = "Poecile montanus"; H = "(Conrad von Baldenstein, 1827)"
    I = $null
    Q = 369436.1354981294; R = 6635294.429866268; S = 25
    Z = "12:52"; AB = "12:52"
}

# Apply the rotation: row4 -> row3 (gets "35"), row5 -> row4, row3 -> row5.
# Clear the old I4 ("35") content first since the new row 4 / row 5 content
# has no value there (old row 5's and row 3's I cells were already empty).
$ws.Range("I4").ClearContents()

Set-RowValues $ws 3 $origRow4
Set-RowValues $ws 4 $origRow5
Set-RowValues $ws 5 $origRow3
